$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 122-123; existing rows 122-137 shift down to 124-139.
$ws.Rows("122:123").Insert()

# New row 122 - Ciruela Angeleno Especial, Provincia de Curicó
$ws.Cells.Item(122,1).Value = 11
$ws.Cells.Item(122,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(122,3).Value = "Bíobío"
$ws.Cells.Item(122,4).Value = 45034
$ws.Cells.Item(122,5).Value = 8
$ws.Cells.Item(122,6).Value = "Fruta"
$ws.Cells.Item(122,7).Value = 100103
$ws.Cells.Item(122,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(122,9).Value = 100103002
$ws.Cells.Item(122,10).Value = "Ciruela"
$ws.Cells.Item(122,11).Value = "Angeleno"
$ws.Cells.Item(122,12).Value = "Especial"
$ws.Cells.Item(122,13).Value = 350
$ws.Cells.Item(122,14).Value = 12000
$ws.Cells.Item(122,15).Value = 13000
$ws.Cells.Item(122,16).Value = 12571
$ws.Cells.Item(122,17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(122,18).Value = "Provincia de Curicó"
$ws.Cells.Item(122,19).Value = 698
$ws.Cells.Item(122,20).Value = 18

# New row 123 - Ciruela Angeleno Primera, Provincia de Curicó
$ws.Cells.Item(123,1).Value = 11
$ws.Cells.Item(123,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(123,3).Value = "Bíobío"
$ws.Cells.Item(123,4).Value = 45034
$ws.Cells.Item(123,5).Value = 8
$ws.Cells.Item(123,6).Value = "Fruta"
$ws.Cells.Item(123,7).Value = 100103
$ws.Cells.Item(123,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(123,9).Value = 100103002
$ws.Cells.Item(123,10).Value = "Ciruela"
$ws.Cells.Item(123,11).Value = "Angeleno"
$ws.Cells.Item(123,12).Value = "Primera"
$ws.Cells.Item(123,13).Value = 300
$ws.Cells.Item(123,14).Value = 10000
$ws.Cells.Item(123,15).Value = 11000
$ws.Cells.Item(123,16).Value = 10500
$ws.Cells.Item(123,17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(123,18).Value = "Provincia de Curicó"
$ws.Cells.Item(123,19).Value = 583
$ws.Cells.Item(123,20).Value = 18
